# Staging.SubSector.xlsx - column/header rework.
#
# The sheet used to carry one column per attribute (A:G) with explicit
# per-column widths. The new layout keeps only column A's width definition
# and rewrites the header row (row 2) with the refreshed attribute names
# (BusinessKey / Code / LongName / SectorBusinessKey / ShortName /
# SubSector_ID / TextDescription), dropping the obsolete "SectorSourceKey"
# column in favor of "SectorBusinessKey".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the per-column width/bestFit formatting that used to live on columns
# B:G - only column A ("width=41, bestFit") should keep an explicit <col>
# definition going forward. This removes the data in B2:G2 too, so we
# re-populate the header row right after with the new attribute names.
$ws.Columns("B:G").Delete()

# Re-write the header row (row 2) with the refreshed attribute names.
$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "LongName"
$ws.Range("D2").Value = "SectorBusinessKey"
$ws.Range("E2").Value = "ShortName"
$ws.Range("F2").Value = "SubSector_ID"
$ws.Range("G2").Value = "TextDescription"

# Reset the active selection back to the top-left cell (the sheet's neutral
# default) instead of the stale "D27" selection left over from editing.
$ws.Range("A1").Select()
